# Auto-generated: apply scheduled market-data refresh to Kraken_Profits workbook.
# For each touched leve row, columns H-N are rewritten to the latest computed
# market snapshot (currentAveragePrice[, NQ, HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]).
# Cells with no resulting value (blank/zero-suppressed columns) are cleared outright.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 367.16666
$ws.Range("I9").Value = 250.8
$ws.Range("K9").Value = 250.8
$ws.Range("M9").Value = -81.80000000000001

# Row 29
$ws.Range("H29").Value = 314.2857
$ws.Range("I29").Value = 200
$ws.Range("K29").Value = 600
$ws.Range("M29").Value = -319

# Row 32
$ws.Range("H32").Value = 6499.3335
$ws.Range("I32").Value = 1665.3334
$ws.Range("J32").Value = 11333.333
$ws.Range("K32").Value = 1665.3334
$ws.Range("L32").Value = 11333.333
$ws.Range("M32").Value = -1339.3334
$ws.Range("N32").Value = -11985.333

# Row 63
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

# Row 66
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

# Row 123
$ws.Range("H123").Value = 99995
$ws.Range("J123").Value = 99995
$ws.Range("L123").Value = 99995
$ws.Range("N123").Value = -109795

# Row 132
$ws.Range("H132").Value = 4132.3335
$ws.Range("I132").Value = 3425.0667
$ws.Range("K132").Value = 10275.2001
$ws.Range("M132").Value = -7745.2001

# Row 137
$ws.Range("H137").Value = 3764
$ws.Range("J137").Value = 3924.5
$ws.Range("L137").Value = 11773.5
$ws.Range("N137").Value = -16873.5

# Row 141
$ws.Range("H141").Value = 15000
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1095.625
$ws.Range("I2").Value = 627.6667
$ws.Range("J2").Value = 2499.5
$ws.Range("K2").Value = 627.6667
$ws.Range("L2").Value = 2499.5
$ws.Range("M2").Value = -514.6667
$ws.Range("N2").Value = -2725.5

# Row 10
$ws.Range("H10").Value = 5000
$ws.Range("I10").Value = 5000
$ws.Range("K10").Value = 5000
$ws.Range("M10").Value = -4830

# Row 11
$ws.Range("H11").Value = 4475.5
$ws.Range("I11").Value = 2451.5
$ws.Range("J11").Value = 6499.5
$ws.Range("K11").Value = 2451.5
$ws.Range("L11").Value = 6499.5
$ws.Range("M11").Value = -2307.5
$ws.Range("N11").Value = -6787.5

# Row 12
$ws.Range("H12").Value = 870
$ws.Range("I12").Value = 300
$ws.Range("J12").Value = 1012.5
$ws.Range("K12").Value = 300
$ws.Range("L12").Value = 1012.5
$ws.Range("M12").Value = -127
$ws.Range("N12").Value = -1358.5

# Row 13
$ws.Range("H13").Value = 23350000
$ws.Range("I13").Value = 25025000
$ws.Range("K13").Value = 25025000
$ws.Range("M13").Value = -25024856

# Row 14
$ws.Range("H14").Value = 500
$ws.Range("I14").Value = 500
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 500
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

# Row 17
$ws.Range("H17").Value = 1009
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

# Row 32
$ws.Range("H32").Value = 1340.6154
$ws.Range("I32").Value = 1434.1666
$ws.Range("K32").Value = 1434.1666
$ws.Range("M32").Value = -1147.1666

# Row 45
$ws.Range("H45").Value = 1812.5
$ws.Range("I45").Value = 1812.5
$ws.Range("K45").Value = 1812.5
$ws.Range("M45").Value = -1435.5

# Row 101
$ws.Range("H101").Value = 12992
$ws.Range("J101").Value = 12992
$ws.Range("L101").Value = 12992
$ws.Range("N101").Value = -19482

# Row 116
$ws.Range("H116").Value = 1095.625
$ws.Range("I116").Value = 627.6667
$ws.Range("J116").Value = 2499.5
$ws.Range("K116").Value = 627.6667
$ws.Range("L116").Value = 2499.5
$ws.Range("M116").Value = 1666.3333
$ws.Range("N116").Value = -7087.5

# Row 122
$ws.Range("H122").Value = 2093.9285
$ws.Range("I122").Value = 2012.1
$ws.Range("K122").Value = 6036.299999999999
$ws.Range("M122").Value = -3586.299999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1095.625
$ws.Range("I3").Value = 627.6667
$ws.Range("J3").Value = 2499.5
$ws.Range("K3").Value = 627.6667
$ws.Range("L3").Value = 2499.5
$ws.Range("M3").Value = -513.6667
$ws.Range("N3").Value = -2727.5

# Row 86
$ws.Range("H86").Value = 1200
$ws.Range("I86").Value = 1300
$ws.Range("J86").Value = 1100
$ws.Range("K86").Value = 1300
$ws.Range("L86").Value = 1100
$ws.Range("M86").Value = -177
$ws.Range("N86").Value = -3346

# Row 89
$ws.Range("H89").Value = 1200
$ws.Range("I89").Value = 1300
$ws.Range("J89").Value = 1100
$ws.Range("K89").Value = 6500
$ws.Range("L89").Value = 5500
$ws.Range("M89").Value = -884
$ws.Range("N89").Value = -16732

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 650
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 50

# Row 31
$ws.Range("H31").Value = 4068.9644
$ws.Range("I31").Value = 1369.6
$ws.Range("J31").Value = 5568.6113
$ws.Range("K31").Value = 1369.6
$ws.Range("L31").Value = 5568.6113
$ws.Range("M31").Value = -1074.6
$ws.Range("N31").Value = -6158.6113

# Row 34
$ws.Range("H34").Value = 4068.9644
$ws.Range("I34").Value = 1369.6
$ws.Range("J34").Value = 5568.6113
$ws.Range("K34").Value = 1369.6
$ws.Range("L34").Value = 5568.6113
$ws.Range("M34").Value = -1167.6
$ws.Range("N34").Value = -5972.6113

# Row 122
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

# Row 132
$ws.Range("H132").Value = 6000
$ws.Range("I132").Value = 6000
$ws.Range("K132").Value = 18000
$ws.Range("M132").Value = -15470

$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 429139.34
$ws.Range("I4").Value = 417246
$ws.Range("K4").Value = 1251738
$ws.Range("M4").Value = -1251626

# Row 15
$ws.Range("H15").Value = 98.75
$ws.Range("J15").Value = 98.333336
$ws.Range("L15").Value = 295.000008
$ws.Range("N15").Value = -575.000008

# Row 23
$ws.Range("H23").Value = 1043.5714
$ws.Range("I23").Value = 1075
$ws.Range("K23").Value = 3225
$ws.Range("M23").Value = -2990

# Row 39
$ws.Range("H39").Value = 4111
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

# Row 60
$ws.Range("H60").Value = 854.44446
$ws.Range("I60").Value = 461.25
$ws.Range("J60").Value = 4000
$ws.Range("K60").Value = 1383.75
$ws.Range("L60").Value = 12000
$ws.Range("M60").Value = -1132.75
$ws.Range("N60").Value = -12502

# Row 61
$ws.Range("H61").Value = 20
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 20
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -490

# Row 98
$ws.Range("H98").Value = 1944
$ws.Range("I98").Value = 1888
$ws.Range("K98").Value = 5664
$ws.Range("M98").Value = -4166

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 3000
$ws.Range("K70").Value = 3000
$ws.Range("M70").Value = -2730

# Row 73
$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 3000
$ws.Range("K73").Value = 3000
$ws.Range("M73").Value = -2064

# Row 102
$ws.Range("H102").Value = 7066.625
$ws.Range("I102").Value = 6423.1665
$ws.Range("K102").Value = 6423.1665
$ws.Range("M102").Value = -4801.1665

# Row 113
$ws.Range("H113").Value = 1600
$ws.Range("I113").Value = 1600
$ws.Range("K113").Value = 1600
$ws.Range("M113").Value = 570

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 4637.375
$ws.Range("I55").Value = 3833.3333
$ws.Range("J55").Value = 5119.8
$ws.Range("K55").Value = 3833.3333
$ws.Range("L55").Value = 5119.8
$ws.Range("M55").Value = -3660.3333
$ws.Range("N55").Value = -5465.8

# Row 74
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

# Row 77
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

# Row 93
$ws.Range("H93").Value = 300
$ws.Range("I93").Value = 300
$ws.Range("K93").Value = 300
$ws.Range("M93").Value = 948

# Row 137
$ws.Range("H137").Value = 67665
$ws.Range("J137").Value = 99995
$ws.Range("L137").Value = 99995
$ws.Range("N137").Value = -110195

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 4058.25
$ws.Range("I122").Value = 4058.25
$ws.Range("K122").Value = 12174.75
$ws.Range("M122").Value = -9724.75
